$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the crypto price ("Price", column D) and volume-change ("Volume(1h)",
# column E) figures in the cryptos list to the values from the latest
# GitHub Actions scrape. A leading apostrophe is used for plain decimal
# price values (e.g. '227.87) so Excel stores them as text instead of
# auto-converting them to numbers, matching the worksheet's existing
# text-formatted price column.

$ws.Cells.Item(2, 4).Value = '34.504.35'
$ws.Cells.Item(2, 5).Value = '  -0.22%  '
$ws.Cells.Item(3, 4).Value = '1.805.99'
$ws.Cells.Item(3, 5).Value = '  -0.74%  '
$ws.Cells.Item(4, 5).Value = '  +0.10%  '
$ws.Cells.Item(5, 4).Value = '''227.87'
$ws.Cells.Item(5, 5).Value = '  -0.06%  '
$ws.Cells.Item(6, 4).Value = '''0.577'
$ws.Cells.Item(6, 5).Value = '  +2.81%  '
$ws.Cells.Item(7, 5).Value = '  +0.07%  '
$ws.Cells.Item(8, 4).Value = '''36.86'
$ws.Cells.Item(8, 5).Value = '  +5.76%  '
$ws.Cells.Item(9, 5).Value = '  -0.35%  '
$ws.Cells.Item(10, 4).Value = '''0.0693'
$ws.Cells.Item(10, 5).Value = '  -0.21%  '
$ws.Cells.Item(11, 5).Value = '  +1.21%  '
$ws.Cells.Item(12, 4).Value = '2.065.73'
$ws.Cells.Item(12, 5).Value = '  -0.87%  '
$ws.Cells.Item(13, 4).Value = '''11.60'
$ws.Cells.Item(13, 5).Value = '  +1.22%  '
$ws.Cells.Item(14, 4).Value = '1.805.20'
$ws.Cells.Item(14, 5).Value = '  -0.93%  '
$ws.Cells.Item(15, 4).Value = '''0.647'
$ws.Cells.Item(15, 5).Value = '  +0.51%  '
$ws.Cells.Item(16, 4).Value = '''4.48'
$ws.Cells.Item(16, 5).Value = '  +3.05%  '
$ws.Cells.Item(17, 4).Value = '34.469.70'
$ws.Cells.Item(17, 5).Value = '  -0.33%  '
$ws.Cells.Item(18, 4).Value = '''70.32'
$ws.Cells.Item(18, 5).Value = '  +1.73%  '
$ws.Cells.Item(19, 4).Value = '''245.60'
$ws.Cells.Item(19, 5).Value = '  -0.58%  '
$ws.Cells.Item(20, 4).Value = '0.0₃0791'
$ws.Cells.Item(20, 5).Value = '  -1.30%  '
$ws.Cells.Item(21, 4).Value = '''11.62'
$ws.Cells.Item(21, 5).Value = '  +0.84%  '
$ws.Cells.Item(22, 5).Value = '  +0.14%  '
$ws.Cells.Item(23, 4).Value = '''4.21'
$ws.Cells.Item(23, 5).Value = '  +0.23%  '
$ws.Cells.Item(24, 4).Value = '''2.18'
$ws.Cells.Item(24, 5).Value = '  +4.27%  '
$ws.Cells.Item(25, 4).Value = '''172.29'
$ws.Cells.Item(25, 5).Value = '  +0.56%  '
$ws.Cells.Item(26, 4).Value = '''8.03'
$ws.Cells.Item(26, 5).Value = '  +8.91%  '
$ws.Cells.Item(28, 5).Value = '  +1.27%  '
$ws.Cells.Item(29, 5).Value = '  +0.00%  '
$ws.Cells.Item(30, 4).Value = '''4.02'
$ws.Cells.Item(30, 5).Value = '  -0.04%  '
$ws.Cells.Item(31, 4).Value = '''3.85'
$ws.Cells.Item(31, 5).Value = '  -0.01%  '
$ws.Cells.Item(32, 5).Value = '  -0.32%  '
$ws.Cells.Item(33, 5).Value = '  -0.44%  '
$ws.Cells.Item(34, 5).Value = '  -1.96%  '
$ws.Cells.Item(35, 4).Value = '1.392.74'
$ws.Cells.Item(35, 5).Value = '  -1.70%  '
$ws.Cells.Item(36, 4).Value = '''0.676'
$ws.Cells.Item(36, 5).Value = '  -0.44%  '
$ws.Cells.Item(37, 4).Value = '''2.46'
$ws.Cells.Item(37, 5).Value = '  -5.68%  '
$ws.Cells.Item(38, 5).Value = '  -0.21%  '
$ws.Cells.Item(39, 5).Value = '  -0.68%  '
$ws.Cells.Item(40, 4).Value = '''83.12'
$ws.Cells.Item(40, 5).Value = '  -3.51%  '
$ws.Cells.Item(41, 4).Value = '''0.963'
$ws.Cells.Item(41, 5).Value = '  +0.84%  '
$ws.Cells.Item(42, 4).Value = '''2.84'
$ws.Cells.Item(42, 5).Value = '  -0.46%  '
$ws.Cells.Item(43, 5).Value = '  +0.46%  '
$ws.Cells.Item(44, 5).Value = '  +7.93%  '
$ws.Cells.Item(45, 4).Value = '''13.50'
$ws.Cells.Item(45, 5).Value = '  -3.35%  '
$ws.Cells.Item(46, 4).Value = '''6.03'
$ws.Cells.Item(46, 5).Value = '  -1.05%  '
$ws.Cells.Item(47, 5).Value = '  -4.54%  '
$ws.Cells.Item(48, 4).Value = '1.967.33'
$ws.Cells.Item(48, 5).Value = '  -0.93%  '
$ws.Cells.Item(49, 4).Value = '''104.58'
$ws.Cells.Item(49, 5).Value = '  -1.26%  '
$ws.Cells.Item(50, 5).Value = '  +0.14%  '
$ws.Cells.Item(51, 5).Value = '  -2.96%  '
